# Atualização de bases das ligas, do dia: 14-06-2024 às 20:31
# Swap the data (columns B through AD, i.e. every field except the
# sequential row index in column A) between pairs of rows so the
# underlying match records trade places while staying on the same
# physical spreadsheet row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($row1, $row2, $firstCol, $lastCol) {
    $values1 = @{}
    $values2 = @{}

    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $values1[$c] = $ws.Cells.Item($row1, $c).Value()
        $values2[$c] = $ws.Cells.Item($row2, $c).Value()
    }

    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($row1, $c).Value = $values2[$c]
        $ws.Cells.Item($row2, $c).Value = $values1[$c]
    }
}

# Column A (index 1) holds the fixed sequential row number and must not
# be touched. Columns B (2) through AD (30) hold the actual record data.
Swap-RowData 28 29 2 30
Swap-RowData 101 102 2 30
Swap-RowData 213 214 2 30
Swap-RowData 263 265 2 30
